$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert a brand-new data row at position 28 (pushes old rows
#    28..46 down to 29..47, dimension becomes A1:T47).
# ------------------------------------------------------------------
$ws.Rows.Item(28).Insert()

# ------------------------------------------------------------------
# 2) Populate the freshly inserted row 28 with the new earthquake
#    record (EventID 5164446). This row is written directly in the
#    final column layout, so no later rotation is applied to it.
# ------------------------------------------------------------------
$ws.Cells.Item(28, 1).Value = 5164446
$ws.Cells.Item(28, 2).Value = 42315.31369212963
$ws.Cells.Item(28, 3).Value = -30.8796
$ws.Cells.Item(28, 4).Value = -71.45189999999999
$ws.Cells.Item(28, 5).Value = 6.8
$ws.Cells.Item(28, 6).Value = "GO03"
$ws.Cells.Item(28, 7).Value = "CO06"
$ws.Cells.Item(28, 8).Value = "CO02"
$ws.Cells.Item(28, 9).Value = "CO03"
$ws.Cells.Item(28, 10).Value = "GO04"
$ws.Cells.Item(28, 11).Value = "TLL"
$ws.Cells.Item(28, 12).Value = "IN40"
$ws.Cells.Item(28, 13).Value = "IN41"
$ws.Cells.Item(28, 14).Value = "2015-11-07T07:31:46"
$ws.Cells.Item(28, 15).Value = "2015-11-07T07:31:49"
$ws.Cells.Item(28, 16).Value = "2015-11-07T07:31:52"
$ws.Cells.Item(28, 17).Value = "2015-11-07T07:31:55"
$ws.Cells.Item(28, 18).Value = "2015-11-07T07:31:55"
$ws.Cells.Item(28, 19).Value = "2015-11-07T07:31:55"
$ws.Cells.Item(28, 20).Value = "2015-11-07T07:31:55"

# ------------------------------------------------------------------
# 3) Re-order columns L:R ("Estacion mas cercana 6/7" plus the five
#    "Inicio..." timestamps) for the header row and for every other
#    data row: the two station codes that used to sit in Q:R move to
#    L:M, and the five timestamps that used to sit in L:P shift right
#    into N:R. This is a rotate-right-by-2 across columns L..R.
# ------------------------------------------------------------------
$ws.Cells.Item(1, 12).Value = "Estaci" + [char]0x00F3 + "n m" + [char]0x00E1 + "s cercana 6"
$ws.Cells.Item(1, 13).Value = "Estaci" + [char]0x00F3 + "n m" + [char]0x00E1 + "s cercana 7"
$ws.Cells.Item(1, 14).Value = "Inicio estaci" + [char]0x00F3 + "n m" + [char]0x00E1 + "s cercana 1"
$ws.Cells.Item(1, 15).Value = "Inicio estaci" + [char]0x00F3 + "n m" + [char]0x00E1 + "s cercana 2"
$ws.Cells.Item(1, 16).Value = "Inicio estaci" + [char]0x00F3 + "n m" + [char]0x00E1 + "s cercana 3"
$ws.Cells.Item(1, 17).Value = "Inicio estaci" + [char]0x00F3 + "n m" + [char]0x00E1 + "s cercana 4"
$ws.Cells.Item(1, 18).Value = "Inicio estaci" + [char]0x00F3 + "n m" + [char]0x00E1 + "s cercana 5"

for ($r = 2; $r -le 47; $r++) {
    if ($r -eq 28) { continue }
    $L = $ws.Cells.Item($r, 12).Text
    $M = $ws.Cells.Item($r, 13).Text
    $N = $ws.Cells.Item($r, 14).Text
    $O = $ws.Cells.Item($r, 15).Text
    $P = $ws.Cells.Item($r, 16).Text
    $Q = $ws.Cells.Item($r, 17).Text
    $R = $ws.Cells.Item($r, 18).Text
    $ws.Cells.Item($r, 12).Value = $Q
    $ws.Cells.Item($r, 13).Value = $R
    $ws.Cells.Item($r, 14).Value = $L
    $ws.Cells.Item($r, 15).Value = $M
    $ws.Cells.Item($r, 16).Value = $N
    $ws.Cells.Item($r, 17).Value = $O
    $ws.Cells.Item($r, 18).Value = $P
}
